$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data (rows 2 and 3), pushing all
# existing data rows down by two (old row 2 -> new row 4, ... old row 120 -> new row 122).
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(3).Insert()

# The Insert() calls copy the header row's (bold/centered) formatting into the
# new rows; reset them back to the plain "Normal" style used by every other
# data row, then re-apply the date number format that column D uses.
$ws.Range("A2:T3").Style = "Normal"
$ws.Range("D2:D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 2
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "Macroferia Regional de Talca"
$ws.Range("C2").Value = "Maule"
$ws.Range("D2").Value = 44860
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100107
$ws.Range("H2").Value = "Otros"
$ws.Range("I2").Value = 100107002
$ws.Range("J2").Value = "Chirimoya"
$ws.Range("K2").Value = "Cultivar IV Región"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 150
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 24000
$ws.Range("Q2").Value = "`$/bandeja 10 kilos"
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 2400
$ws.Range("T2").Value = 10

# New row 3
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "Macroferia Regional de Talca"
$ws.Range("C3").Value = "Maule"
$ws.Range("D3").Value = 44860
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107002
$ws.Range("J3").Value = "Chirimoya"
$ws.Range("K3").Value = "Cultivar IV Región"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 22000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 22000
$ws.Range("Q3").Value = "`$/bandeja 10 kilos"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 2200
$ws.Range("T3").Value = 10
